$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.259.33"
$ws.Range("E2").Value = "  -5.78%  "

$ws.Range("D3").Value = "2.458.73"
$ws.Range("E3").Value = "  -8.32%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'540.16"
$ws.Range("E5").Value = "  -2.77%  "

$ws.Range("D6").Value = "'146.31"
$ws.Range("E6").Value = "  -7.20%  "

$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("D8").Value = "'0.575"
$ws.Range("E8").Value = "  -2.39%  "

$ws.Range("D9").Value = "2.477.84"
$ws.Range("E9").Value = "  -7.68%  "

$ws.Range("D10").Value = "'0.0998"
$ws.Range("E10").Value = "  -5.38%  "

$ws.Range("E11").Value = "  -1.51%  "

$ws.Range("D12").Value = "'5.53"
$ws.Range("E12").Value = "  +0.47%  "

$ws.Range("D13").Value = "'0.352"
$ws.Range("E13").Value = "  -3.89%  "

$ws.Range("D14").Value = "2.890.49"
$ws.Range("E14").Value = "  -8.45%  "

$ws.Range("D15").Value = "'24.24"
$ws.Range("E15").Value = "  -7.66%  "

$ws.Range("D16").Value = "59.225.42"
$ws.Range("E16").Value = "  -5.69%  "

$ws.Range("D17").Value = "'0.0000139"
$ws.Range("E17").Value = "  -5.28%  "

$ws.Range("D18").Value = "2.499.81"
$ws.Range("E18").Value = "  -6.85%  "

$ws.Range("D19").Value = "'11.21"
$ws.Range("E19").Value = "  -5.38%  "

$ws.Range("D20").Value = "'4.37"
$ws.Range("E20").Value = "  -4.92%  "

$ws.Range("D21").Value = "'324.34"
$ws.Range("E21").Value = "  -5.99%  "

$ws.Range("D22").Value = "'0.968"
$ws.Range("E22").Value = "  -3.22%  "

$ws.Range("D23").Value = "'5.75"
$ws.Range("E23").Value = "  -7.50%  "

$ws.Range("D24").Value = "'60.86"
$ws.Range("E24").Value = "  -3.86%  "

$ws.Range("D25").Value = "'0.452"
$ws.Range("E25").Value = "  -11.20%  "

$ws.Range("D26").Value = "'0.161"
$ws.Range("E26").Value = "  -5.06%  "

$ws.Range("D27").Value = "'0.977"
$ws.Range("E27").Value = "  -2.26%  "

$ws.Range("D28").Value = "'7.77"
$ws.Range("E28").Value = "  -4.81%  "

$ws.Range("D29").Value = "'6.83"
$ws.Range("E29").Value = "  -6.16%  "

$ws.Range("E30").Value = "  -9.30%  "

$ws.Range("D31").Value = "'1.83"
$ws.Range("E31").Value = "  -4.99%  "

$ws.Range("D32").Value = "0.0₃0777"
$ws.Range("E32").Value = "  -8.80%  "

$ws.Range("E33").Value = "  -0.17%  "

$ws.Range("D34").Value = "'157.78"
$ws.Range("E34").Value = "  -3.31%  "

$ws.Range("D35").Value = "'1.42"
$ws.Range("E35").Value = "  -2.86%  "

$ws.Range("D36").Value = "'18.69"
$ws.Range("E36").Value = "  -3.89%  "

$ws.Range("D37").Value = "'4.48"
$ws.Range("E37").Value = "  -7.97%  "

$ws.Range("D38").Value = "'1.72"
$ws.Range("E38").Value = "  -3.84%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'5.85"
$ws.Range("E39").Value = "  -5.20%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'313.57"
$ws.Range("E40").Value = "  -7.99%  "

$ws.Range("D41").Value = "'36.44"
$ws.Range("E41").Value = "  -4.82%  "

$ws.Range("D42").Value = "'3.73"
$ws.Range("E42").Value = "  -6.56%  "

$ws.Range("D43").Value = "'0.833"
$ws.Range("E43").Value = "  -10.25%  "

$ws.Range("E44").Value = "  -0.27%  "

$ws.Range("D45").Value = "'0.599"
$ws.Range("E45").Value = "  -2.67%  "

$ws.Range("D46").Value = "'10.73"
$ws.Range("E46").Value = "  -2.45%  "

$ws.Range("D47").Value = "'0.0531"
$ws.Range("E47").Value = "  -4.71%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'125.05"
$ws.Range("E48").Value = "  -4.44%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.0936"
$ws.Range("E49").Value = "  -3.48%  "

$ws.Range("D50").Value = "'0.0231"
$ws.Range("E50").Value = "  -4.21%  "

$ws.Range("D51").Value = "'18.45"
$ws.Range("E51").Value = "  -8.45%  "
